# Fix the wording in Questions.xlsx ("easy" sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("easy")

$ws.Range("A19").Value = "If you choose a square at random from the electrode array, what is the probability of hitting any part of the neuron? "
$ws.Range("A20").Value = "If you choose a square at random from the electrode array, what is the probability of hitting the axon hillock?"
$ws.Range("A21").Value = "If you choose a square at random from the electrode array, what is the probability of hitting any part of the axon?"
$ws.Range("A22").Value = "If you choose a square at random from the electrode array, what is the probability of hitting any part of the dendrite?"
$ws.Range("A24").Value = "Awesome! We will find our battleship in no time. Lets further learn about sets."
$ws.Range("A31").Value = "Are P(D) and P(AH) mutually exclusive?"

# Update the view state to match the final selection/scroll position
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 14
$ws.Range("A28").Select()
